$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.874.54"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "3.212.31"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "604.91"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.19%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "153.42"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.210.18"
$ws.Range("E8").Value = "  +1.24%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.532"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("E11").Value = "  -1.59%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.509"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("E13").Value = "  +0.78%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "38.90"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.23%  "

$ws.Range("D15").Value = "3.737.99"
$ws.Range("E15").Value = "  +1.32%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.47"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.28%  "

$ws.Range("D17").Value = "66.033.76"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").Value = "3.198.48"
$ws.Range("E18").Value = "  +0.91%  "

$ws.Range("E19").Value = "  +0.03%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "509.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("E21").Value = "  +4.39%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.737"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("E23").Value = "  +0.41%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.07"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.44%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "85.11"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +2.78%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.94%  "

$ws.Range("E30").Value = "  +3.52%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.79"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.27%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "28.12"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +0.41%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "55.01"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.44%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0903"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.88%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "477.90"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0419"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("E40").Value = "  -5.22%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.31%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.297"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.60%  "

$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("D44").Value = "2.948.48"
$ws.Range("E44").Value = "  -3.71%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("D46").Value = "0.0₃0638"
$ws.Range("E46").Value = "  +5.24%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "28.54"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  +0.94%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.68%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "121.28"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.08%  "
